# Updates cryptos list figures (price + 1h volume change) per the
# "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Helper: write a value that must stay a text cell (even when it looks
# numeric, e.g. "242.34") without leaving the cell style changed -
# mark it as Text, assign, then drop the format override again.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$ws.Range("D2").Value = "29.332.22"
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("D3").Value = "1.878.18"
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("E4").Value = "  +0.17%  "
Set-TextValue $ws.Range("D5") "0.7108"
$ws.Range("E5").Value = "  -0.47%  "
Set-TextValue $ws.Range("D6") "242.34"
$ws.Range("E6").Value = "  +0.21%  "
Set-TextValue $ws.Range("D8") "0.08001"
Set-TextValue $ws.Range("D9") "0.3163"
$ws.Range("E9").Value = "  +1.71%  "
Set-TextValue $ws.Range("D10") "24.99"
$ws.Range("E10").Value = "  -0.48%  "
Set-TextValue $ws.Range("D11") "0.08306"
$ws.Range("E11").Value = "  -1.63%  "
$ws.Range("D12").Value = "1.879.86"
$ws.Range("E12").Value = "  +0.00%  "
Set-TextValue $ws.Range("D13") "5.251"
$ws.Range("E13").Value = "  -0.10%  "
Set-TextValue $ws.Range("D14") "94.36"
$ws.Range("E14").Value = "  +3.51%  "
Set-TextValue $ws.Range("D15") "0.7139"
$ws.Range("E15").Value = "  +0.18%  "
Set-TextValue $ws.Range("D16") "6.368"
$ws.Range("E16").Value = "  +4.66%  "
Set-TextValue $ws.Range("D17") "0.000008510"
$ws.Range("E17").Value = "  +3.45%  "
$ws.Range("D18").Value = "29.339.78"
Set-TextValue $ws.Range("D19") "243.48"
$ws.Range("E19").Value = "  +1.15%  "
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue $ws.Range("D20") "13.28"
$ws.Range("E20").Value = "  +0.41%  "
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "2.128.92"
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("E24").Value = "  +0.17%  "
Set-TextValue $ws.Range("D25") "0.1560"
$ws.Range("E25").Value = "  -2.23%  "
Set-TextValue $ws.Range("D26") "9.069"
$ws.Range("E26").Value = "  +0.21%  "
Set-TextValue $ws.Range("D27") "162.75"
$ws.Range("E27").Value = "  -0.26%  "
Set-TextValue $ws.Range("D28") "18.55"
$ws.Range("E28").Value = "  +0.09%  "
Set-TextValue $ws.Range("D29") "1.506"
Set-TextValue $ws.Range("D30") "4.419"
$ws.Range("E30").Value = "  +0.02%  "
Set-TextValue $ws.Range("D31") "4.323"
$ws.Range("E31").Value = "  -0.20%  "
Set-TextValue $ws.Range("D32") "1.193"
$ws.Range("E32").Value = "  -7.00%  "
Set-TextValue $ws.Range("D33") "0.05388"
$ws.Range("E33").Value = "  +1.45%  "
Set-TextValue $ws.Range("D34") "1.937"
$ws.Range("E34").Value = "  -0.02%  "
Set-TextValue $ws.Range("D35") "0.7717"
$ws.Range("E35").Value = "  +4.27%  "
Set-TextValue $ws.Range("D36") "1.182"
$ws.Range("E36").Value = "  +0.30%  "
Set-TextValue $ws.Range("D37") "2.682"
$ws.Range("E37").Value = "  -0.60%  "
$ws.Range("E38").Value = "  +0.90%  "
$ws.Range("D39").Value = "1.261.49"
$ws.Range("E39").Value = "  +2.36%  "
Set-TextValue $ws.Range("D40") "2.753"
$ws.Range("E40").Value = "  +0.84%  "
Set-TextValue $ws.Range("D41") "6.486"
$ws.Range("E41").Value = "  -0.53%  "
Set-TextValue $ws.Range("D42") "113.07"
$ws.Range("E42").Value = "  +2.33%  "
Set-TextValue $ws.Range("D43") "0.9051"
$ws.Range("E43").Value = "  +1.54%  "
Set-TextValue $ws.Range("D44") "74.16"
$ws.Range("E44").Value = "  +1.57%  "
$ws.Range("E45").Value = "  +7.66%  "
$ws.Range("E46").Value = "  +0.14%  "
$ws.Range("D47").Value = "2.029.93"
$ws.Range("E47").Value = "  +0.38%  "
Set-TextValue $ws.Range("D48") "0.5228"
$ws.Range("E48").Value = "  +0.30%  "
Set-TextValue $ws.Range("D49") "1.801"
$ws.Range("E49").Value = "  -0.49%  "
Set-TextValue $ws.Range("D50") "9.456"
$ws.Range("E50").Value = "  +0.12%  "
Set-TextValue $ws.Range("D51") "0.4366"
$ws.Range("E51").Value = "  +1.19%  "
